# Apply "Add data for 2022-07-16" update:
# - Rename sheet/title date from 2022-07-07 to 2022-07-08
# - Update the "July (through 07-07)" label to "July (through 07-08)"
# - Update July row (row 8) and Total row (row 9) values for each year column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab title
$ws.Name = "Through 2022-07-08"

# Update the "July (through 07-07)" label in column A, row 8
$ws.Range("A8").Value = "July (through 07-08)"

# Update July row (row 8) values, columns B..I
$julyOld = @{ "B8" = 9; "C8" = 11; "D8" = 12; "E8" = 23; "F8" = 11; "G8" = 26; "H8" = 41; "I8" = 37 }
foreach ($addr in $julyOld.Keys) {
    $ws.Range($addr).Value = $julyOld[$addr]
}

# Update Total row (row 9) values, columns B..I
$totalNew = @{ "B9" = 134; "C9" = 259; "D9" = 402; "E9" = 376; "F9" = 262; "G9" = 498; "H9" = 801; "I9" = 843 }
foreach ($addr in $totalNew.Keys) {
    $ws.Range($addr).Value = $totalNew[$addr]
}
